# Applies crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.424.03"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "3.681.87"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'686.95"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("D6").Value = "'159.63"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.38%  "

$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("D10").Value = "'7.09"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("D11").Value = "'0.434"
$ws.Range("E11").Value = "  -3.46%  "

$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").Value = "4.304.27"

$ws.Range("D14").Value = "'32.24"
$ws.Range("E14").Value = "  -3.33%  "

$ws.Range("D15").Value = "69.439.20"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "3.669.14"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("E17").Value = "  +2.04%  "

$ws.Range("D18").Value = "'15.81"
$ws.Range("E18").Value = "  -2.82%  "

$ws.Range("D19").Value = "'6.38"
$ws.Range("E19").Value = "  -3.49%  "

$ws.Range("D20").Value = "'470.55"
$ws.Range("E20").Value = "  -2.48%  "

$ws.Range("D21").Value = "'9.93"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").Value = "'0.649"
$ws.Range("E22").Value = "  -1.89%  "

$ws.Range("D23").Value = "'79.59"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "3.828.89"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("E26").Value = "  -2.41%  "

$ws.Range("D27").Value = "'10.96"
$ws.Range("E27").Value = "  -5.06%  "

$ws.Range("D28").Value = "'9.18"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("E30").Value = "  -4.60%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.58"
$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.99"
$ws.Range("E32").Value = "  -5.32%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'26.83"

$ws.Range("D35").Value = "3.656.03"
$ws.Range("E35").Value = "  +0.33%  "

$ws.Range("E36").Value = "  -1.97%  "

$ws.Range("D37").Value = "'8.20"
$ws.Range("E37").Value = "  -3.32%  "

$ws.Range("D38").Value = "'6.13"
$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  +2.06%  "

$ws.Range("E41").Value = "  -4.41%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("E43").Value = "  -1.50%  "

$ws.Range("D44").Value = "'165.61"
$ws.Range("E44").Value = "  +5.60%  "

$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").Value = "'0.000281"
$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("E48").Value = "  +5.83%  "

$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").Value = "'27.59"
$ws.Range("E50").Value = "  -2.51%  "

$ws.Range("E51").Value = "  -3.05%  "
